# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Each hashtable maps a row number to its new value; the same set of events
# is refreshed with updated "want to go" counts on both sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds "想去人数"
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    4  = 408
    5  = 1003
    6  = 5594
    7  = 505
    8  = 708
    11 = 81
    17 = 1883
    19 = 943
    23 = 563
    24 = 162
    25 = 1058
    28 = 3023
    30 = 107
    31 = 68
    32 = 130
    34 = 413
    39 = 299
    44 = 72
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Sheet "全部类型" (all types) - same events appear again, column F holds "想去人数"
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    5  = 1003
    7  = 5594
    8  = 505
    9  = 708
    16 = 81
    23 = 1883
    25 = 943
    29 = 563
    30 = 162
    31 = 1058
    32 = 3023
    34 = 107
    35 = 68
    36 = 130
    38 = 413
    42 = 299
    46 = 72
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
